# ModifyBuiltinProperties.pptx - remove the Aspose.Slides "evaluation"
# watermark textboxes left over from the trial build, and tidy up the
# stray placeholder left behind on the second slide.

$p = $ppt.ActivePresentation

# --- Slide 1 -----------------------------------------------------------
# The evaluation watermark ("Evaluation only. / Created with ... /
# Copyright ...") is a standalone TextBox shape - delete it outright.
$s1 = $p.Slides.Item(1)
for ($i = $s1.Shapes.Count; $i -ge 1; $i--) {
    $shp = $s1.Shapes.Item($i)
    if ($shp.Name -eq "TextBox") {
        $shp.Delete()
    }
}

# --- Slide 2 -------------------------------------------------------------
# Same watermark shape exists here too, but instead of removing it we
# clear its text (the three evaluation paragraphs) and shrink it down to
# a tiny leftover placeholder, then drop a small new empty textbox next
# to it (what's left after manually wiping out the watermark by hand).
$s2 = $p.Slides.Item(2)
$watermark = $null
for ($i = 1; $i -le $s2.Shapes.Count; $i++) {
    $shp = $s2.Shapes.Item($i)
    if ($shp.Name -eq "TextBox") {
        $watermark = $shp
    }
}

if ($watermark -ne $null) {
    $watermark.TextFrame.TextRange.Text = ""
    $watermark.Left = 4479841 / 12700
    $watermark.Top = 3051623 / 12700
    $watermark.Width = 184730 / 12700
    $watermark.Height = 754694 / 12700
}

$newBox = $s2.Shapes.AddTextbox(1, 3995936 / 12700, 2924944 / 12700, 184731 / 12700, 369332 / 12700)
$newBox.TextFrame.WordWrap = 0
$newBox.TextFrame.AutoSize = 1
$newBox.Fill.Visible = 0
